$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row
$ws.Range("F2").Value = "SchemeStartDate"
$ws.Range("H2").Value = "SchemeEndTime"
$ws.Range("J2").Value = "ReservationID"

# Add new ReservationID values for reserved (non-Available) rows
$ws.Range("J4").Value = 9
$ws.Range("J5").Value = 10
$ws.Range("J11").Value = 11
$ws.Range("J12").Value = 12
$ws.Range("J14").Value = 8
$ws.Range("J15").Value = 8
$ws.Range("J16").Value = 8
$ws.Range("J17").Value = 8
$ws.Range("J21").Value = 13
$ws.Range("J22").Value = 14
$ws.Range("J24").Value = 15
$ws.Range("J25").Value = 16
$ws.Range("J26").Value = 17
$ws.Range("J27").Value = 18

# Update selected cell / view
[void]$ws.Range("I30").Select()
